# Apply the diff: divide values in D8:D13 and H8:H13 by 10 (fixes a units
# error in the measurement table), and update the active selection on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Divide raw measured values by 10 (D8:D13 and H8:H13).
foreach ($row in 8..13) {
    $dCell = $ws.Cells.Item($row, 4)
    $dCell.Value = $dCell.Value() / 10   # column D

    $hCell = $ws.Cells.Item($row, 8)
    $hCell.Value = $hCell.Value() / 10   # column H
}

# Update the selection on the sheet to a single cell H14 (was D6:J13 with active cell D6).
$ws.Activate()
$ws.Range("H14").Select()

$wb.Save()
